# Refresh the crypto price/volume snapshot (scheduled GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new "Price" values are numeric strings with a trailing
# zero (e.g. "1.000", "9.160", "11.60"). Left alone, Excel auto-converts a
# numeric-looking .Value assignment to a real Number and drops that trailing
# zero. Pre-formatting those specific cells as Text keeps the literal string
# exactly as published.
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '28.132.08'
$ws.Range('E2').Value = '  -0.37%  '

# Row 3
$ws.Range('D3').Value = '1.829.83'
$ws.Range('E3').Value = '  +1.50%  '

# Row 4
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.29%  '

# Row 5
$ws.Range('D5').Value = '311.07'
$ws.Range('E5').Value = '  -1.05%  '

# Row 6
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.34%  '

# Row 7
$ws.Range('D7').Value = '0.5118'
$ws.Range('E7').Value = '  -2.58%  '

# Row 8
$ws.Range('D8').Value = '0.3967'
$ws.Range('E8').Value = '  +3.79%  '

# Row 9
$ws.Range('D9').Value = '0.09978'
$ws.Range('E9').Value = '  +24.41%  '

# Row 10
$ws.Range('D10').Value = '1.112'
$ws.Range('E10').Value = '  +0.96%  '

# Row 11
$ws.Range('D11').Value = '41.01'
$ws.Range('E11').Value = '  -0.71%  '

# Row 12
$ws.Range('D12').Value = '6.483'
$ws.Range('E12').Value = '  +2.67%  '

# Row 13
$ws.Range('D13').Value = '1.001'
$ws.Range('E13').Value = '  -0.26%  '

# Row 14
$ws.Range('D14').Value = '20.67'
$ws.Range('E14').Value = '  +0.40%  '

# Row 15
$ws.Range('D15').Value = '7.416'
$ws.Range('E15').Value = '  +1.39%  '

# Row 16
$ws.Range('D16').Value = '1.820.68'
$ws.Range('E16').Value = '  +0.69%  '

# Row 17
$ws.Range('D17').Value = '0.00001135'
$ws.Range('E17').Value = '  +3.67%  '

# Row 18
$ws.Range('D18').Value = '94.15'
$ws.Range('E18').Value = '  +2.16%  '

# Row 19
$ws.Range('D19').Value = '0.06626'
$ws.Range('E19').Value = '  +0.29%  '

# Row 20
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '0.9994'
$ws.Range('E20').Value = '  -0.38%  '

# Row 21
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = '17.38'
$ws.Range('E21').Value = '  +0.08%  '

# Row 22
$ws.Range('D22').Value = '6.056'
$ws.Range('E22').Value = '  +1.41%  '

# Row 23
$ws.Range('D23').Value = '28.208.73'
$ws.Range('E23').Value = '  -0.29%  '

# Row 24
$ws.Range('D24').Value = '11.19'
$ws.Range('E24').Value = '  +0.45%  '

# Row 25
$ws.Range('D25').Value = '2.249'
$ws.Range('E25').Value = '  -1.31%  '

# Row 26
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '20.87'
$ws.Range('E26').Value = '  +1.92%  '

# Row 27
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.466'
$ws.Range('E27').Value = '  +4.60%  '

# Row 28
$ws.Range('D28').Value = '158.24'
$ws.Range('E28').Value = '  -1.55%  '

# Row 29
$ws.Range('D29').Value = '2.035.60'
$ws.Range('E29').Value = '  +1.25%  '

# Row 30
$ws.Range('D30').Value = '128.85'
$ws.Range('E30').Value = '  +4.48%  '

# Row 31
$ws.Range('E31').Value = '  +1.40%  '

# Row 32
$ws.Range('D32').Value = '1.064'
$ws.Range('E32').Value = '  +0.60%  '

# Row 33
$ws.Range('D33').Value = '5.644'
$ws.Range('E33').Value = '  +1.59%  '

# Row 34
$ws.Range('D34').Value = '3.641'
$ws.Range('E34').Value = '  -1.30%  '

# Row 35
$ws.Range('D35').Value = '0.06903'
$ws.Range('E35').Value = '  -4.56%  '

# Row 36
$ws.Range('D36').Value = '9.160'
$ws.Range('E36').Value = '  +6.41%  '

# Row 37
$ws.Range('D37').Value = '0.02344'
$ws.Range('E37').Value = '  +1.42%  '

# Row 38
$ws.Range('D38').Value = '0.2173'
$ws.Range('E38').Value = '  +1.20%  '

# Row 39
$ws.Range('D39').Value = '11.60'
$ws.Range('E39').Value = '  -6.44%  '

# Row 40
$ws.Range('D40').Value = '5.034'
$ws.Range('E40').Value = '  -1.65%  '

# Row 41
$ws.Range('D41').Value = '0.6289'
$ws.Range('E41').Value = '  +1.47%  '

# Row 42
$ws.Range('D42').Value = '0.9994'
$ws.Range('E42').Value = '  -0.25%  '

# Row 43
$ws.Range('E43').Value = '  -1.18%  '

# Row 44
$ws.Range('D44').Value = '13.37'
$ws.Range('E44').Value = '  +0.64%  '

# Row 45
$ws.Range('D45').Value = '0.6008'
$ws.Range('E45').Value = '  -0.11%  '

# Row 46
$ws.Range('D46').Value = '1.292'
$ws.Range('E46').Value = '  -5.85%  '

# Row 47
$ws.Range('D47').Value = '3.711'
$ws.Range('E47').Value = '  -1.55%  '

# Row 48
$ws.Range('D48').Value = '125.83'
$ws.Range('E48').Value = '  -1.08%  '

# Row 49
$ws.Range('D49').Value = '1.996'
$ws.Range('E49').Value = '  +3.60%  '

# Row 50
$ws.Range('D50').Value = '1.191'
$ws.Range('E50').Value = '  -2.56%  '

# Row 51
$ws.Range('D51').Value = '0.06793'
$ws.Range('E51').Value = '  -0.17%  '
